# The source workbook tracks weekly "Betarraga" price observations in
# Sheet1, one row per observation, ordered from oldest-appended to
# newest-appended going down the sheet (rows 2..593). This commit adds a
# new weekly observation.
#
# The new record is inserted at row 552 (pushing the former rows
# 552..593 down to 553..594, so the sheet grows from A1:R593 to
# A1:R594). The new row at 552 is a near-duplicate of the (now shifted)
# row 553 -- same market/category/quality/volume/prices -- except for
# the date (column D) and the weighted-average price (column M), which
# carry the new observation's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 552
$lastCol = 18   # column R

# Push rows 552..593 down to 553..594, creating a fresh blank row 552.
$ws.Rows.Item($newRow).Insert()

# Populate the new row 552 by copying the record that is now sitting at
# row 553 (the original row-552 record, shifted down) ...
for ($col = 1; $col -le $lastCol; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $ws.Cells.Item($newRow + 1, $col).Value2
}

# ... then overwrite the two cells that actually differ for the new
# observation: the date and the weighted average price.
$ws.Cells.Item($newRow, 4).Value = 45223   # column D (Fecha)
$ws.Cells.Item($newRow, 13).Value = 627    # column M (Precio promedio ponderado)
